$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.813.87"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.632.46"
$ws.Range("E3").Value = "  +3.76%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'515.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Value = "'143.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'0.565"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "2.658.97"
$ws.Range("E9").Value = "  +4.71%  "
$ws.Range("D10").Value = "'6.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("D12").Value = "'0.335"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").Value = "3.099.47"
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D15").Value = "58.785.98"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "'20.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "2.650.11"
$ws.Range("E18").Value = "  +4.31%  "
$ws.Range("D19").Value = "'347.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.95%  "
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "'10.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("D22").Value = "'6.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.39%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'61.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "'0.162"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").Value = "0.0₃0799"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").Value = "'7.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  +7.40%  "
$ws.Range("D32").Value = "'18.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("E33").Value = "  +2.37%  "
$ws.Range("D34").Value = "'149.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "'0.967"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.36%  "
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("E40").Value = "  +5.21%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Value = "'277.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").Value = "'0.611"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("D44").Value = "'0.995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "'0.0984"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").Value = "'19.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.78%  "
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").Value = "'10.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").Value = "1.980.93"
$ws.Range("E50").Value = "  +4.29%  "
$ws.Range("D51").Value = "'4.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.26%  "
